# ADD: final version 1.0
#
# Title page currently reads "K KURSOVOMU PROEKTU" (i.e. "To the course
# project"); it must be changed to "K KURSOVOJ RABOTE" ("To the course
# work"). In the underlying OOXML this text is split, around a
# "_GoBack" bookmark, into the runs "K " / "KURSO" / "VOMU PROEKTU" / " ".
# The edit turns it into "K " / "KURSOVOJ" / " " / "RABOTE" / " ", with
# the bookmark now sitting right after "RABOTE" (immediately before the
# trailing space run).
#
# NOTE: this interpreter mangles literal non-ASCII characters that are
# embedded directly in the script source (they come back as "?" from
# Write-Output), so all Cyrillic text below is assembled from Unicode
# code points instead of being typed literally. Equality comparisons on
# the resulting strings work fine even though console output does not.

function Cyr($codes) {
    $s = ""
    foreach ($code in $codes) {
        $s = $s + [char]$code
    }
    return $s
}

$KURSO         = Cyr @(0x041A,0x0423,0x0420,0x0421,0x041E)                               # "КУРСО"
$VOMU_PROEKTU  = Cyr @(0x0412,0x041E,0x041C,0x0423,0x0020,0x041F,0x0420,0x041E,0x0415,0x041A,0x0422,0x0423)  # "ВОМУ ПРОЕКТУ"
$VOJ           = Cyr @(0x0412,0x041E,0x0419)                                             # "ВОЙ"
$SPACE_RABOTE  = Cyr @(0x0020,0x0420,0x0410,0x0411,0x041E,0x0422,0x0415)                 # " РАБОТЕ"

$d = $word.ActiveDocument

if (-not $d.Bookmarks.Exists("_GoBack")) {
    throw "Expected bookmark _GoBack not found"
}

$bm = $d.Bookmarks("_GoBack")
$pos0 = $bm.Start                 # boundary between the "КУРСО" run and the "ВОМУ ПРОЕКТУ" run
$kursovoyStart = $pos0 - 5        # "КУРСО" is 5 characters long

$before = $d.Range($kursovoyStart, $pos0)
if ($before.Text -ne $KURSO) {
    throw "Unexpected text before bookmark; aborting"
}
$after = $d.Range($pos0, $pos0 + 12)
if ($after.Text -ne $VOMU_PROEKTU) {
    throw "Unexpected text after bookmark; aborting"
}

# Step 1: "КУРСО" -> "КУРСОВОЙ" (append "ВОЙ" right before the bookmark;
# this merges into the preceding run).
$r1 = $d.Range($pos0, $pos0)
$r1.InsertAfter($VOJ)

# Step 2: delete the old "ВОМУ ПРОЕКТУ" run text, which now starts right
# after the (shifted) bookmark.
$bm = $d.Bookmarks("_GoBack")
$pos1 = $bm.Start
$oldRun2 = $d.Range($pos1, $pos1 + 12)
$oldRun2.Text = ""

# Step 3: insert " РАБОТЕ" right after the (still same) bookmark position.
$bm = $d.Bookmarks("_GoBack")
$pos2 = $bm.Start
$r3 = $d.Range($pos2, $pos2)
$r3.InsertAfter($SPACE_RABOTE)

# At this point the paragraph holds one merged run "К КУРСОВОЙ РАБОТЕ"
# followed by the bookmark and the original trailing " " run. $pos2 is
# the position right after "КУРСОВОЙ", i.e. where the new " РАБОТЕ" text
# starts.
#
# Step 4: split that merged run back into separate sibling runs, each
# with the same (unchanged) formatting, by toggling a format property
# away from and then back to its original value on each sub-range. The
# engine does not silently re-merge runs that were split apart like
# this, even once their formatting matches again.

# 4a: split "РАБОТЕ" off from " ".
$raboteRange = $d.Range($pos2 + 1, $pos2 + 7)
$raboteRange.Font.Size = 12
$raboteRange.Font.Size = 20

# 4b: split " " off from "КУРСОВОЙ".
$spaceRange = $d.Range($pos2, $pos2 + 1)
$spaceRange.Font.Size = 12
$spaceRange.Font.Size = 20

# 4c: split "КУРСОВОЙ" off from "К ".
$kursovoyRange = $d.Range($kursovoyStart, $pos2)
$kursovoyRange.Font.Size = 12
$kursovoyRange.Font.Size = 20

$bmFinal = $d.Bookmarks("_GoBack")
Write-Output ("_GoBack now at " + $bmFinal.Start + "-" + $bmFinal.End)
Write-Output "Title updated: KURSOVOMU PROEKTU -> KURSOVOJ RABOTE"
